# Add new columns I (I0) and J (IF) to Sheet1, matching headers and borders
# of the existing header row, and fill values for rows 2-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 / J1 ---
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Match formatting of the other header cells (bold font, thin border,
# centered horizontal/top vertical alignment) as used by H1.
$headerRange = $ws.Range($ws.Cells.Item(1, 9), $ws.Cells.Item(1, 10))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows 2-49: column I ("I0") and column J ("IF") ---
$rowValues = @{
    2 = @(8, 8)
    3 = @(9, 9)
    4 = @(9, 9)
    5 = @(7, 7)
    6 = @(8, 8)
    7 = @(8, 8)
    8 = @(7, 7)
    9 = @(8, 8)
    10 = @(10, 10)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(9, 9)
    20 = @(8, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(8, 8)
    25 = @(9, 9)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(8, 8)
    29 = @(9, 9)
    30 = @(8, 8)
    31 = @(8, 8)
    32 = @(8, 8)
    33 = @(8, 8)
    34 = @(7, 7)
    35 = @(8, 8)
    36 = @(7, 8)
    37 = @(7, 7)
    38 = @(7, 7)
    39 = @(9, 9)
    40 = @(8, 8)
    41 = @(7, 7)
    42 = @(7, 7)
    43 = @(8, 8)
    44 = @(8, 8)
    45 = @(8, 8)
    46 = @(8, 8)
    47 = @(8, 8)
    48 = @(8, 8)
    49 = @(6, 6)
}

foreach ($r in $rowValues.Keys) {
    $pair = $rowValues[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}

Write-Output "I0 and IF columns populated"
